$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 468.13635
$ws.Range("I32").Value2 = 418.4
$ws.Range("J32").Value2 = 482.7647
$ws.Range("K32").Value2 = 418.4
$ws.Range("L32").Value2 = 482.7647
$ws.Range("M32").Value2 = -92.39999999999998
$ws.Range("N32").Value2 = -1134.7647

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value2 = 519.5
$ws.Range("I41").Value2 = 106.14286
$ws.Range("J41").Value2 = 1098.2
$ws.Range("K41").Value2 = 106.14286
$ws.Range("L41").Value2 = 1098.2
$ws.Range("M41").Value2 = 333.85714
$ws.Range("N41").Value2 = -1978.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value2 = 1590617.8
$ws.Range("I98").Value2 = 1794858.6
$ws.Range("J98").Value2 = 7750.75
$ws.Range("K98").Value2 = 1794858.6
$ws.Range("L98").Value2 = 7750.75
$ws.Range("M98").Value2 = -1793360.6
$ws.Range("N98").Value2 = -10746.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value2 = 1590617.8
$ws.Range("I122").Value2 = 1794858.6
$ws.Range("J122").Value2 = 7750.75
$ws.Range("K122").Value2 = 5384575.800000001
$ws.Range("L122").Value2 = 23252.25
$ws.Range("M122").Value2 = -5382125.800000001
$ws.Range("N122").Value2 = -28152.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value2 = 32748.889
$ws.Range("I123").Value2 = 0
$ws.Range("J123").Value2 = 32748.889
$ws.Range("K123").Value2 = 0
$ws.Range("L123").Value2 = 32748.889
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value2 = -42548.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value2 = 2148.889
$ws.Range("I131").Value2 = 691.4286
$ws.Range("J131").Value2 = 7250
$ws.Range("K131").Value2 = 2074.2858
$ws.Range("L131").Value2 = 21750
$ws.Range("M131").Value2 = 2965.7142
$ws.Range("N131").Value2 = -31830

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 3060.2368
$ws.Range("I138").Value2 = 1102.3914
$ws.Range("J138").Value2 = 6062.2666
$ws.Range("K138").Value2 = 3307.1742
$ws.Range("L138").Value2 = 18186.7998
$ws.Range("M138").Value2 = 1832.8258
$ws.Range("N138").Value2 = -28466.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1998.625
$ws.Range("I61").Value2 = 2012.6786
$ws.Range("J61").Value2 = 1965.8334
$ws.Range("K61").Value2 = 2012.6786
$ws.Range("L61").Value2 = 1965.8334
$ws.Range("M61").Value2 = -1800.6786
$ws.Range("N61").Value2 = -2389.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value2 = 5500.857
$ws.Range("I88").Value2 = 1835.3334
$ws.Range("J88").Value2 = 8250
$ws.Range("K88").Value2 = 1835.3334
$ws.Range("L88").Value2 = 8250
$ws.Range("M88").Value2 = -1429.3334
$ws.Range("N88").Value2 = -9062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value2 = 5500.857
$ws.Range("I91").Value2 = 1835.3334
$ws.Range("J91").Value2 = 8250
$ws.Range("K91").Value2 = 1835.3334
$ws.Range("L91").Value2 = 8250
$ws.Range("M91").Value2 = -431.3334
$ws.Range("N91").Value2 = -11058

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 1998.625
$ws.Range("I136").Value2 = 2012.6786
$ws.Range("J136").Value2 = 1965.8334
$ws.Range("K136").Value2 = 6038.0358
$ws.Range("L136").Value2 = 5897.5002
$ws.Range("M136").Value2 = -3488.0358
$ws.Range("N136").Value2 = -10997.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 1842.1489
$ws.Range("I134").Value2 = 1747.9756
$ws.Range("K134").Value2 = 5243.9268
$ws.Range("M134").Value2 = -2708.9268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1893.0714
$ws.Range("I16").Value2 = 1345.35
$ws.Range("J16").Value2 = 3262.375
$ws.Range("K16").Value2 = 1345.35
$ws.Range("L16").Value2 = 3262.375
$ws.Range("M16").Value2 = -1058.35
$ws.Range("N16").Value2 = -3836.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 23810512
$ws.Range("I58").Value2 = 47619650
$ws.Range("J58").Value2 = 1373.5238
$ws.Range("K58").Value2 = 47619650
$ws.Range("L58").Value2 = 1373.5238
$ws.Range("M58").Value2 = -47619447
$ws.Range("N58").Value2 = -1779.5238

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value2 = 1893.0714
$ws.Range("I113").Value2 = 1345.35
$ws.Range("J113").Value2 = 3262.375
$ws.Range("K113").Value2 = 1345.35
$ws.Range("L113").Value2 = 3262.375
$ws.Range("M113").Value2 = 824.6500000000001
$ws.Range("N113").Value2 = -7602.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 23810512
$ws.Range("I136").Value2 = 47619650
$ws.Range("J136").Value2 = 1373.5238
$ws.Range("K136").Value2 = 142858950
$ws.Range("L136").Value2 = 4120.5714
$ws.Range("M136").Value2 = -142856400
$ws.Range("N136").Value2 = -9220.571400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value2 = 0
$ws.Range("I66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("K66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 881.56384
$ws.Range("J113").Value2 = 963.2655999999999
$ws.Range("L113").Value2 = 2889.7968
$ws.Range("N113").Value2 = -7229.7968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value2 = 0
$ws.Range("I114").Value2 = 0
$ws.Range("J114").Value2 = 0
$ws.Range("K114").Value2 = 0
$ws.Range("L114").Value2 = 0
$ws.Range("M114").ClearContents()
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 488.5
$ws.Range("I122").Value2 = 372.30768
$ws.Range("J122").Value2 = 1999
$ws.Range("K122").Value2 = 3350.76912
$ws.Range("L122").Value2 = 17991
$ws.Range("M122").Value2 = -900.7691199999999
$ws.Range("N122").Value2 = -22891

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value2 = 1034.3636
$ws.Range("I132").Value2 = 799.75
$ws.Range("J132").Value2 = 1660
$ws.Range("K132").Value2 = 7197.75
$ws.Range("L132").Value2 = 14940
$ws.Range("M132").Value2 = -4667.75
$ws.Range("N132").Value2 = -20000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value2 = 4500
$ws.Range("I6").Value2 = 1000
$ws.Range("J6").Value2 = 8000
$ws.Range("K6").Value2 = 1000
$ws.Range("L6").Value2 = 8000
$ws.Range("M6").Value2 = -887
$ws.Range("N6").Value2 = -8226

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H16").Value2 = 4500
$ws.Range("I16").Value2 = 1000
$ws.Range("J16").Value2 = 8000
$ws.Range("K16").Value2 = 1000
$ws.Range("L16").Value2 = 8000
$ws.Range("M16").Value2 = -750
$ws.Range("N16").Value2 = -8500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value2 = 7864.6
$ws.Range("I18").Value2 = 7005
$ws.Range("J18").Value2 = 7960.1113
$ws.Range("K18").Value2 = 7005
$ws.Range("L18").Value2 = 7960.1113
$ws.Range("M18").Value2 = -6712
$ws.Range("N18").Value2 = -8546.1113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 488.47726
$ws.Range("I22").Value2 = 445.32144
$ws.Range("J22").Value2 = 564
$ws.Range("K22").Value2 = 445.32144
$ws.Range("L22").Value2 = 564
$ws.Range("M22").Value2 = -150.32144
$ws.Range("N22").Value2 = -1154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value2 = 488.47726
$ws.Range("I27").Value2 = 445.32144
$ws.Range("J27").Value2 = 564
$ws.Range("K27").Value2 = 445.32144
$ws.Range("L27").Value2 = 564
$ws.Range("M27").Value2 = -338.32144
$ws.Range("N27").Value2 = -778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value2 = 553
$ws.Range("I30").Value2 = 553
$ws.Range("K30").Value2 = 553
$ws.Range("M30").Value2 = -445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 3176135
$ws.Range("I132").Value2 = 3922461
$ws.Range("K132").Value2 = 11767383
$ws.Range("M132").Value2 = -11764853

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 51890.85
$ws.Range("I136").Value2 = 144785.42
$ws.Range("J136").Value2 = 1870.6923
$ws.Range("K136").Value2 = 434356.26
$ws.Range("L136").Value2 = 5612.0769
$ws.Range("M136").Value2 = -431806.26
$ws.Range("N136").Value2 = -10712.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value2 = 3000
$ws.Range("J20").Value2 = 3000
$ws.Range("L20").Value2 = 3000
$ws.Range("N20").Value2 = -3480

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value2 = 423.33334
$ws.Range("I113").Value2 = 335
$ws.Range("J113").Value2 = 600
$ws.Range("K113").Value2 = 1005
$ws.Range("L113").Value2 = 1800
$ws.Range("M113").Value2 = 1165
$ws.Range("N113").Value2 = -6140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 1242.9062
$ws.Range("I122").Value2 = 1039.3704
$ws.Range("J122").Value2 = 2342
$ws.Range("K122").Value2 = 3118.1112
$ws.Range("L122").Value2 = 7026
$ws.Range("M122").Value2 = -668.1112000000003
$ws.Range("N122").Value2 = -11926
